$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HANA flow")

# Swap the "Project Name" lookup value between rows: rows 6, 10, 11, 16 move
# from "CDS Item Operations Extended" (TWC4618) to "CDS Item Operations"
# (TWC3149), while row 18 moves the other way. Columns I/J/K/L recompute
# automatically off these via existing sheet formulas.
$ws.Range("H6").Value = "CDS Item Operations"
$ws.Range("H10").Value = "CDS Item Operations"
$ws.Range("H11").Value = "CDS Item Operations"
$ws.Range("H16").Value = "CDS Item Operations"
$ws.Range("H18").Value = "CDS Item Operations Extended"

# Reflect the edit location as the sheet's active selection.
$ws.Activate() | Out-Null
$ws.Range("H6").Select() | Out-Null
